$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.206.33"
$ws.Range("E2").Value = "  -2.45%  "
$ws.Range("D3").Value = "1.671.71"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "217.21"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "0.5120"
$ws.Range("E6").Value = "  -1.32%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "0.2649"
$ws.Range("E8").Value = "  +1.51%  "
$ws.Range("D9").Value = "0.06374"
$ws.Range("E9").Value = "  +2.98%  "
$ws.Range("D10").Value = "21.52"
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("D11").Value = "0.07385"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.669.10"
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.533"
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").Value = "0.5811"
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").Value = "0.000008606"
$ws.Range("E15").Value = "  +5.26%  "
$ws.Range("D16").Value = "64.31"
$ws.Range("E16").Value = "  -2.18%  "
$ws.Range("D17").Value = "26.189.12"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("D18").Value = "4.927"
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "190.09"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "6.199"
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "144.89"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("D25").Value = "7.626"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").Value = "0.1177"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").Value = "15.65"
$ws.Range("E27").Value = "  +2.23%  "
$ws.Range("D28").Value = "0.06087"
$ws.Range("E28").Value = "  +5.97%  "
$ws.Range("D29").Value = "1.294"
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("D31").Value = "3.519"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "3.523"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "1.636"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D34").Value = "1.014"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("D35").Value = "0.6045"
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("D36").Value = "2.380"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("D37").Value = "2.664"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("D38").Value = "6.171"
$ws.Range("E38").Value = "  +5.15%  "
$ws.Range("D39").Value = "0.01609"
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "1.076.03"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("D41").Value = "0.8657"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "100.83"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("D44").Value = "1.822.51"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("E45").Value = "  +7.77%  "
$ws.Range("D46").Value = "56.14"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").Value = "8.079"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").Value = "0.05212"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "0.4297"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "1.434"
$ws.Range("E51").Value = "  -1.88%  "
